$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column C (Control Name) - target stored width is 42.6640625 characters.
# The COM ColumnWidth setter here quantizes to 1/6-character steps, so 41.83 is the
# input that lands closest to the target stored width.
$ws.Columns.Item(3).ColumnWidth = 41.83

# Fill in the previously-empty Test Type / Cloud Resource / Category / Responsibility /
# Validation Steps / USNORTHCOM Validated columns (H:M) for every data row with "-".
$ws.Range("H2:M57").Value = "-"

# Row 14 (AC-02_IAM_User_No_Policies_Check) got specific tagging values instead of "-".
$ws.Range("H14").Value = "AWS Managed"
$ws.Range("I14").Value = "IAM"
$ws.Range("J14").Value = "Not Inherited"
$ws.Range("K14").Value = "USNORTHCOM"
